$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the very start of the document (inside
# paragraph 1). The edit relocates it to the middle of the new paragraph about the
# COIL-20 dataset, so drop the old one first to avoid a duplicate-name clash.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p2  = $d.Paragraphs.Item(2)
$p15 = $d.Paragraphs.Item(15)
$region = $d.Range($p2.Range.Start, $p15.Range.End)

$payload = @'
<w:p><w:r><w:t>ECE 5332-011: Deep Learning for Medical Signal/Image data</w:t></w:r></w:p><w:p><w:r><w:t>Alexandre Soares</w:t></w:r><w:r><w:t xml:space="preserve"> da Silva, R11485685</w:t></w:r></w:p><w:p><w:r><w:t>Rishi</w:t></w:r></w:p><w:p><w:r><w:t>Jud</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Changing the radius of </w:t></w:r><w:r><w:t xml:space="preserve">the high/low pass filters increases or decreases their modifying power in the images. </w:t></w:r><w:r><w:t>Form smaller radius sizes,</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>COMMENT ABOUT CONVOLUTION AND SPACIAL FILTERING</w:t></w:r></w:p><w:p><w:r><w:t>COMMENT ABOUT CLUSTERING</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>3. Simple Image classification problem</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The proposed problem statement required </w:t></w:r><w:r><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t>use</w:t></w:r><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:r><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t xml:space="preserve">SURF (Speeded-up Robust Features) </w:t></w:r><w:r><w:t xml:space="preserve">algorithm to match and classify </w:t></w:r><w:r><w:t xml:space="preserve">10 3D-rendered objects from 24 different perspectives. The </w:t></w:r><w:r><w:t xml:space="preserve">assignment’s </w:t></w:r><w:r><w:t xml:space="preserve">goal was </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>evaluate</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the accuracy of </w:t></w:r><w:r><w:t>such classification engine with the use of different numbers of features and varying dataset sizes for training and testing.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">The dataset used was the </w:t></w:r><w:r><w:t>Columbia University Image Library</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:t>COIL-20</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> dataset</w:t></w:r><w:r><w:t xml:space="preserve">, with selection of 10 objects </w:t></w:r><w:r><w:t>in 24 different perspectives</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p/>
'@

$xml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $payload + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$region.InsertXML($xml)

Write-Output "Paragraph count after edit: $($d.Paragraphs.Count)"
